# Update countries & provincias Spain
#
# The source feed refreshed: a handful of countries' case counts moved,
# which re-ranks some neighbouring rows in the (already case-count sorted)
# table, and the "last updated" footer timestamp advanced from 03:35 to
# 04:05.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Footer / header timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 04:05"

# --- Row 58: Australia ---
$ws.Range("B58").Value = 7081
$ws.Range("C58").Value = 2
$ws.Range("D58").Value = 6470
$ws.Range("E58").Value = 511

# --- Row 67: Bolivia ---
$ws.Range("B67").Value = 4919
$ws.Range("C67").Value = 438
$ws.Range("D67").Value = 553
$ws.Range("E67").Value = 4167
$ws.Range("G67").Value = 10
$ws.Range("H67").Value = 199

# --- Rows 82-85: Guatemala overtakes Croacia, Costa de Marfil y Tayikistan ---
$ws.Range("A82").Value = "Guatemala"
$ws.Range("B82").Value = 2265
$ws.Range("C82").Value = 132
$ws.Range("D82").Value = 159
$ws.Range("E82").Value = 2061
$ws.Range("G82").Value = 2
$ws.Range("H82").Value = 45

$ws.Range("A83").Value = "Croacia"
$ws.Range("B83").Value = 2234
$ws.Range("D83").Value = 1978
$ws.Range("E83").Value = 160
$ws.Range("H83").Value = 96

$ws.Range("A84").Value = "Costa de Marfil"
$ws.Range("B84").Value = 2231
$ws.Range("D84").Value = 1083
$ws.Range("E84").Value = 1119
$ws.Range("H84").Value = 29

$ws.Range("A85").Value = "Tayikistan"
$ws.Range("B85").Value = 2140
$ws.Range("D85").Value = 470
$ws.Range("E85").Value = 1629
$ws.Range("H85").Value = 41

# --- Row 163: Benin ---
$ws.Range("D163").Value = 57
$ws.Range("E163").Value = 71

# --- Rows 197-198: Santa Lucia overtakes Belice ---
$ws.Range("A197").Value = "Santa Lucia"
$ws.Range("D197").Value = 18
$ws.Range("H197").Value = 0

$ws.Range("A198").Value = "Belice"
$ws.Range("D198").Value = 16
$ws.Range("H198").Value = 2

# --- Rows 209-210: Montserrat overtakes Groenlandia ---
$ws.Range("A209").Value = "Montserrat"
$ws.Range("D209").Value = 10
$ws.Range("H209").Value = 1

$ws.Range("A210").Value = "Groenlandia"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0
